$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the table with a new "2021" column (M), mirroring the existing
# year columns (D:L). Copy column L's formatting (borders, number format,
# etc.) into column M for the header block (rows 3-7), then overwrite the
# values for the new year.
$ws.Range("L3:L7").Copy() | Out-Null
$ws.Range("M3:M7").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("M4").Value = 2021
$ws.Range("M5").Value = 98
$ws.Range("M6").Value = 97
$ws.Range("M7").Value = 96

# Reset the view: scroll back to the top-left and select A1 (the sheet was
# left scrolled/selected at B1/N13 before the edit).
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("A1").Select() | Out-Null
